$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 39, pushing the existing rows 39-45 down to 40-46
$ws.Rows.Item(39).Insert()

# Populate the new row 39 with the new weekly data point (same dimension values
# as the row that follows it, except for the fields that actually changed).
$ws.Range("A39").Value = 11
$ws.Range("B39").Value = "Vega Monumental Concepción"
$ws.Range("C39").Value = "Bíobío"
$ws.Range("D39").Value = 45135
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 100112022
$ws.Range("G39").Value = "Arveja Verde"
$ws.Range("H39").Value = "Perfection"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 18000
$ws.Range("L39").Value = 20000
$ws.Range("M39").Value = 19000
$ws.Range("N39").Value = "$/malla 25 kilos"
$ws.Range("O39").Value = "Provincia de Limarí"
$ws.Range("P39").Value = 760
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
